$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.206.80"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "1.602.19"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("D4").Value = "0.9992"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "0.9993"
$ws.Range("D6").Value = "303.33"
$ws.Range("E6").Value = "  +0.68%  "
$ws.Range("D7").Value = "0.3781"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "51.64"
$ws.Range("E8").Value = "  +3.44%  "
$ws.Range("D9").Value = "0.3619"
$ws.Range("E9").Value = "  -0.66%  "
$ws.Range("D10").Value = "1.267"
$ws.Range("E10").Value = "  +0.77%  "
$ws.Range("D11").Value = "0.9998"
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("D12").Value = "0.08123"
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("D13").Value = "22.67"
$ws.Range("E13").Value = "  -1.54%  "
$ws.Range("D14").Value = "6.607"
$ws.Range("E14").Value = "  -0.24%  "
$ws.Range("D15").Value = "7.410"
$ws.Range("E15").Value = "  +0.20%  "
$ws.Range("E16").Value = "  -0.24%  "
$ws.Range("D17").Value = "1.605.10"
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("D18").Value = "93.56"
$ws.Range("E18").Value = "  +2.04%  "
$ws.Range("D19").Value = "0.06861"
$ws.Range("E19").Value = "  +0.42%  "
$ws.Range("D20").Value = "18.08"
$ws.Range("E20").Value = "  -0.94%  "
$ws.Range("D21").Value = "6.533"
$ws.Range("E21").Value = "  -0.52%  "
$ws.Range("D22").Value = "0.9988"
$ws.Range("E22").Value = "  -0.25%  "
$ws.Range("D23").Value = "12.98"
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("D24").Value = "23.198.05"
$ws.Range("E24").Value = "  +0.27%  "
$ws.Range("B25").Value = "LidoDAOToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D25").Value = "3.037"
$ws.Range("E25").Value = "  +10.41%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "2.390"
$ws.Range("E26").Value = "  +2.18%  "
$ws.Range("D27").Value = "21.19"
$ws.Range("E27").Value = "  +0.44%  "
$ws.Range("D28").Value = "149.92"
$ws.Range("E28").Value = "  -0.40%  "
$ws.Range("D29").Value = "5.231"
$ws.Range("E29").Value = "  -0.99%  "
$ws.Range("D30").Value = "133.91"
$ws.Range("E30").Value = "  +1.59%  "
$ws.Range("D31").Value = "2.425"
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("D32").Value = "6.830"
$ws.Range("E32").Value = "  -0.22%  "
$ws.Range("D33").Value = "1.783.76"
$ws.Range("E33").Value = "  +0.39%  "
$ws.Range("D34").Value = "0.9866"
$ws.Range("E34").Value = "  +4.48%  "
$ws.Range("D35").Value = "0.07574"
$ws.Range("E35").Value = "  -1.64%  "
$ws.Range("D36").Value = "10.38"
$ws.Range("E36").Value = "  +2.98%  "
$ws.Range("D37").Value = "0.02728"
$ws.Range("E37").Value = "  -0.67%  "
$ws.Range("D38").Value = "6.166"
$ws.Range("E38").Value = "  -1.38%  "
$ws.Range("D39").Value = "0.2509"
$ws.Range("E39").Value = "  -1.23%  "
$ws.Range("D40").Value = "0.08793"
$ws.Range("E40").Value = "  -1.29%  "
$ws.Range("D41").Value = "0.7130"
$ws.Range("E41").Value = "  +0.32%  "
$ws.Range("E42").Value = "  -2.13%  "
$ws.Range("E43").Value = "  -2.87%  "
$ws.Range("D44").Value = "15.62"
$ws.Range("E44").Value = "  +0.42%  "
$ws.Range("D45").Value = "0.6573"
$ws.Range("E45").Value = "  -0.70%  "
$ws.Range("D46").Value = "2.316"
$ws.Range("E46").Value = "  +0.74%  "
$ws.Range("D47").Value = "4.016"
$ws.Range("E47").Value = "  +1.05%  "
$ws.Range("D48").Value = "132.26"
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("D49").Value = "0.07959"
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("D50").Value = "1.210"
$ws.Range("E50").Value = "  -0.11%  "
$ws.Range("E51").Value = "  +3.29%  "
